$d = $word.ActiveDocument

# Remove the three paragraphs that follow the "Consolidação das Leis do
# Trabalho..." reference: the blank separator paragraph, the
# "Ver no Jupiter..." paragraph, and the "© 2020 ..." footer paragraph.
# Deleting from the end keeps earlier paragraph indices valid.

$footer = $d.Content.Find
$footer.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($footer.Found) {
    $p = $footer.Parent.Paragraphs(1)
    $p.Range.Delete()
}

$jupiter = $d.Content.Find
$jupiter.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($jupiter.Found) {
    $p = $jupiter.Parent.Paragraphs(1)
    $p.Range.Delete()
}

$consolid = $d.Content.Find
$consolid.Execute("Consolidação das Leis do Trabalho. Rio De Janeiro: Forense, 1994.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($consolid.Found) {
    $consolidPara = $consolid.Parent.Paragraphs(1)
    $nextPara = $consolidPara.Next()
    $nextPara.Range.Delete()
}
